$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.346.25"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.848.70"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.34"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2735"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06297"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").Value = "1.830.99"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07448"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.945"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.02"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6218"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "30.299.03"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.16"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007307"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.931"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.885"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.203"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.879"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1027"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.076"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.818"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04865"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.148"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7132"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.704"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01894"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.661"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8839"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "105.85"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.919"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.547"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4037"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.34"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1203"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.625"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.28"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05517"
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.357"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3664"
$ws.Range("E51").Value = "  -1.19%  "
